$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.397.76'
$ws.Range("E2").Value = '  +1.60%  '

$ws.Range("D3").Value = '1.866.34'
$ws.Range("E3").Value = '  +1.42%  '

$ws.Range("E4").Value = '  +0.26%  '

$ws.Range("D5").Value = "'239.40"

$ws.Range("E6").Value = '  +0.90%  '

$ws.Range("E7").Value = '  +0.24%  '

$ws.Range("D8").Value = "'42.93"
$ws.Range("E8").Value = '  +8.27%  '

$ws.Range("D9").Value = "'0.332"
$ws.Range("E9").Value = '  +0.53%  '

$ws.Range("E10").Value = '  +1.56%  '

$ws.Range("D11").Value = "'0.0989"
$ws.Range("E11").Value = '  +0.87%  '

$ws.Range("E12").Value = '  +1.31%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.901.34'
$ws.Range("E13").Value = '  +3.53%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = "'11.58"
$ws.Range("E14").Value = '  +1.65%  '

$ws.Range("D15").Value = "'0.683"
$ws.Range("E15").Value = '  +1.58%  '

$ws.Range("D16").Value = "'4.74"
$ws.Range("E16").Value = '  +2.06%  '

$ws.Range("D17").Value = '35.352.37'
$ws.Range("E17").Value = '  +1.40%  '

$ws.Range("D18").Value = "'70.26"
$ws.Range("E18").Value = '  +0.64%  '

$ws.Range("E19").Value = '  +1.66%  '

$ws.Range("D20").Value = "'241.35"
$ws.Range("E20").Value = '  +0.27%  '

$ws.Range("E21").Value = '  +0.92%  '

$ws.Range("D22").Value = "'4.76"
$ws.Range("E22").Value = '  +1.51%  '

$ws.Range("E23").Value = '  +0.32%  '

$ws.Range("D24").Value = "'2.25"
$ws.Range("E24").Value = '  -1.26%  '

$ws.Range("D25").Value = "'169.72"
$ws.Range("E25").Value = '  -1.02%  '

$ws.Range("D26").Value = "'1.90"
$ws.Range("E26").Value = '  +24.89%  '

$ws.Range("E27").Value = '  +5.47%  '

$ws.Range("D28").Value = "'17.78"
$ws.Range("E28").Value = '  +1.88%  '

$ws.Range("E29").Value = '  +1.44%  '

$ws.Range("D30").Value = "'0.0563"
$ws.Range("E30").Value = '  +2.19%  '

$ws.Range("E31").Value = '  +0.26%  '

$ws.Range("D33").Value = "'1.83"
$ws.Range("E33").Value = '  +26.57%  '

$ws.Range("D34").Value = "'4.06"
$ws.Range("E34").Value = '  +2.69%  '

$ws.Range("E35").Value = '  +9.38%  '

$ws.Range("D36").Value = "'0.821"
$ws.Range("E36").Value = '  +18.53%  '

$ws.Range("E37").Value = '  +6.13%  '

$ws.Range("E38").Value = '  +3.60%  '

$ws.Range("E39").Value = '  +4.76%  '

$ws.Range("D40").Value = "'91.01"
$ws.Range("E40").Value = '  +0.79%  '

$ws.Range("D41").Value = '1.349.54'
$ws.Range("E41").Value = '  +0.36%  '

$ws.Range("E42").Value = '  +15.71%  '

$ws.Range("D43").Value = "'15.20"
$ws.Range("E43").Value = '  +1.72%  '

$ws.Range("B44").Value = 'Gas'
$ws.Range("C44").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D44").Value = "'13.22"
$ws.Range("E44").Value = '  +56.87%  '

$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = "'2.36"
$ws.Range("E45").Value = '  +2.57%  '

$ws.Range("D46").Value = "'2.41"
$ws.Range("E46").Value = '  +0.60%  '

$ws.Range("D47").Value = "'6.66"
$ws.Range("E47").Value = '  +6.37%  '

$ws.Range("E48").Value = '  -1.14%  '

$ws.Range("D49").Value = '2.050.81'
$ws.Range("E49").Value = '  +1.53%  '

$ws.Range("D50").Value = "'0.0688"
$ws.Range("E50").Value = '  +3.35%  '

$ws.Range("D51").Value = "'3.41"
$ws.Range("E51").Value = '  -1.07%  '

# Reset style on cells that were force-quoted as text, to avoid stray quote-prefix formatting
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
